$p = $ppt.ActivePresentation

$NL = [char]13

function Set-ParagraphLines($shape, [string[]]$lines) {
    $tr = $shape.TextFrame.TextRange
    # Replace with an unrelated placeholder first so the COM host's
    # text-diffing doesn't try to preserve/split runs against the old
    # text (it matches common substrings and fragments <a:r> runs).
    $tr.Text = "."
    $tr.Text = [string]::Join($NL, $lines)
}

# --- Slide 2: Introduction ---
$s2 = $p.Slides.Item(2)
$lines2 = @(
    "Recurrent neural networks have been state-of-the-art for sequence modeling and transduction tasks .",
    "However , they suffer from sequential computation , limiting parallelization and preventing longer sequence lengths .",
    "Attention mechanisms have become integral to sequence modeling , but they are typically used with recurrent networks ."
)
Set-ParagraphLines $s2.Shapes.Item(2) $lines2

# --- Slide 3: Literature Review ---
$s3 = $p.Slides.Item(3)
$lines3 = @(
    "nan ; Recent advances in computational efficiency include factorization tricks and conditional computation .",
    "However , fundamental constraint of sequential computation remains .",
    "Attention mechanisms have become integral part of compelling sequence modeling and transduction models in various tasks , allowing modeling of dependencies without regard to their distance in input or output sequences .",
    "In all but few cases , however , such attention mechanisms are used in conjunction with recurrent network .",
    "; nan"
)
Set-ParagraphLines $s3.Shapes.Item(2) $lines3

# --- Slide 4: Methodology ---
$s4 = $p.Slides.Item(4)
$lines4 = @(
    "Transformer follows encoder-decoder structure with stacked self-attention and point-wise , fully connected layers for both encoder and decoder .",
    "Attention is computed using scaled dot-product attention , where weights are computed by compatibility function of query with corresponding key .",
    "Multi-head attention consists of several attention layers running in parallel .",
    "Positional encodings are added to input embeddings to inject information about relative or absolute position of tokens in sequence ."
)
Set-ParagraphLines $s4.Shapes.Item(2) $lines4

# --- Slide 5: Results ---
$s5 = $p.Slides.Item(5)
$lines5 = @(
    "Transformer achieves state-of-the-art results on two machine translation tasks , outperforming existing best results , including ensembles , by over 2 BLEU on WMT 2014 English-to-German translation task .",
    "Transformer establishes new single-model state of art score of 41.8 after training for 3.5 days on eight GPUs .",
    "It also generalizes well to other tasks , such as English constituency parsing ."
)
Set-ParagraphLines $s5.Shapes.Item(2) $lines5

# --- Slide 6: Conclusion ---
$s6 = $p.Slides.Item(6)
$lines6 = @(
    "Transformer introduces new architecture for sequence transduction tasks , which is based solely on attention mechanisms and dispenses with recurrence and convolution entirely .",
    "This architecture allows for significantly more parallelization and can reach new state of art in translation quality after being trained for as little as twelve hours on eight P100 GPUs .",
    "Transformer also generalizes well to other tasks such as English constituency parsing ."
)
Set-ParagraphLines $s6.Shapes.Item(2) $lines6
